$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new regression-test row 21 (3004 epic regression case).
# Write column C (the shared long path string) before column B so the
# shared-strings table gets the new unique strings in the same order as
# the target workbook (path string first, label string second).
$ws.Range("A21").Value = "test"
$ws.Range("C21").Value = "\Testdata\Non_Oncology\DataFiles\Protocol_Page\DownloadProtocol\Protocol_PRISMA_Data.xlsx"
$ws.Range("B21").Value = "download_protocol_prisma"

$ws.Range("A21").Select()
